$d = $word.ActiveDocument
$sel = $word.Selection
$sel.EndKey(6, 0)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Cody Rhodes</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Trine University</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>SE353: Software Engineering</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak /><w:t>SE353 Project Phase 1</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Abstract</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Project Party is our temporary name for a party game with asymmetric online multiplayer, custom user-made characters, and minigames. The game features one host player and everyone else as contestants. The goal is to be the first contestant to reach the end of the board while the host hinders your progress. The game is going to be written in C# using the Unity Engine.  </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Project Description</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>1.</w:t></w:r><w:r><w:tab /><w:t>Main Goals:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>To create an easily accessible party game with a fun gameplay loop that is quick to learn.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>2.</w:t></w:r><w:r><w:tab /><w:t>Main Functionality and Characteristics:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>a.</w:t></w:r><w:r><w:tab /><w:t>Asymmetric Multiplayer</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>The game shall be based around Asymmetric multiplayer. One player shall be the Host, and the rest of the players will be Contestants. The Host shall have actions that can either aid or harm the Contestants’ progression towards the end of the board.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>b.</w:t></w:r><w:r><w:tab /><w:t>Character Creation:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart" /><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t>.</w:t></w:r><w:r><w:tab /><w:t>Players shall be able to choose from several premade contestants or build their own contestants to use in the game. This generation will be done by selecting a model and then allocating stats.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>ii.</w:t></w:r><w:r><w:tab /><w:t xml:space="preserve">The stats for a contestant are health and speed. Stats shall be allocated using a point buy system. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>c.</w:t></w:r><w:r><w:tab /><w:t>Level Design:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The game shall have different levels called “Boards” and all shall have their own pool of minigames. </w:t></w:r><w:proofErr w:type="gramStart" /><w:r><w:t>All of</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:t xml:space="preserve"> the minigames on per board shall be themed to match the board. For example, a beach themed board may have a beach ball spiking minigame. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>d.</w:t></w:r><w:r><w:tab /><w:t>User Friendly GUI:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak /><w:t xml:space="preserve">Navigating the home menu and the interactions in the game shall be simple and follow accepted standards. Meaning the Menus shall not be cluttered, and the user can easily identify where they are or where they need to go. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>e.</w:t></w:r><w:r><w:tab /><w:t>The Loop</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart" /><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t>.</w:t></w:r><w:r><w:tab /><w:t xml:space="preserve">Every player shall take an action on their turn controlling their contestant. Turn order will be determined by speed. If a tie for speed occurs the game will have the players roll a die. The higher die wins and gets to go first. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>ii.</w:t></w:r><w:r><w:tab /><w:t>Contestants shall compete against each other to reach the end of the board. Upon reaching the end of the board that player shall win the game. The other contestants shall continue playing to determine further placements.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>iii.</w:t></w:r><w:r><w:tab /><w:t>The competitions shall be a variety of minigames that are tailored to the boards.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>iv.</w:t></w:r><w:r><w:tab /><w:t xml:space="preserve">If a contestant reaches zero health, they will become unconscious for a turn before returning to 1 health. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>v.</w:t></w:r><w:r><w:tab /><w:t xml:space="preserve">After contestants </w:t></w:r><w:proofErr w:type="gramStart" /><w:r><w:t>reach</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:t xml:space="preserve"> the end of the board, they shall become specters that can interact with the remaining players in a similar fashion to the host but not as powerful. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>vi.</w:t></w:r><w:r><w:tab /><w:t xml:space="preserve">The host shall have a variety of abilities that they are able to use when it is their turn. Their abilities shall be benevolent or harmful to all players so that the host cannot pick a favorite player and help them win. For example, one of their abilities may be to bring a thunderstorm in and all contestants will be struck by lightning reducing their movement temporarily and health permanently. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>3.</w:t></w:r><w:r><w:tab /><w:t>Intended Users and Key Usability Goals</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>a.</w:t></w:r><w:r><w:tab /><w:t>An intended user of this game are people who like the Mario Party series of games but either only play on PC or just want to change things up. Another intended user is a person of any age who wants a simple game to play with their friends.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>4.</w:t></w:r><w:r><w:tab /><w:t>Intended Technology</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>a.</w:t></w:r><w:r><w:tab /><w:t xml:space="preserve">For this project we are planning on using Unity as our Game Engine and the Unity C# API. To handle the </w:t></w:r><w:proofErr w:type="gramStart" /><w:r><w:t>multiplayer</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:t xml:space="preserve"> we are still looking for solutions, but we are wanting to use a peer to peer connection system where all information is routed through the lobby host’s console. For creating models, we will be using Blender. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>5.</w:t></w:r><w:r><w:tab /><w:t>Potential Challenges</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak /><w:t>a.</w:t></w:r><w:r><w:tab /><w:t xml:space="preserve">Neither of us have experience writing </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:t>netcode</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:t xml:space="preserve"> so that will be the main challenge for us. Another challenge is not being creative enough to make unique minigames for every map. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>6.</w:t></w:r><w:r><w:tab /><w:t>Potential for Further Development</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>a.</w:t></w:r><w:r><w:tab /><w:t>To further develop the game more boards could be added, and more game modes could be added to add more variety.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>7.</w:t></w:r><w:r><w:tab /><w:t>Professional Growth</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>a.</w:t></w:r><w:r><w:tab /><w:t xml:space="preserve">Both of us are looking to get careers in Game Design and/or Development. This project shall look impressive on both of our resumes. The skills we shall learn from this project such as the team management and networking will be very beneficial. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>8.</w:t></w:r><w:r><w:tab /><w:t>About the Team Members</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Mason Bledsoe - Extended Reality Major from Shelburn, Indiana. Experience with C++, C# Python and the Unity Engine. Will be focused on the Game Design, character creation, and networking.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Cody Rhodes – Extended Reality Student from Monee, Illinois. Experience with C++, C#, Java, Python, and working on a team. Will be focusing on developing the gameplay mechanics.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Market Potential</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Our game is going to have a lot of market potential. It is going to be relatively inexpensive and will be widely available. It can be listed on Steam as well as any other online game marketplace. We are going to make it easy for the average user to run so their computer should not be a limiting factor when they decide if they want the game or not. We also have a good idea that people will be interested in our game. Party games have always been popular. Thinking back to their time as a child everyone can remember a game or two that they played with a group of friends. They have not always been video games but even since party video games were </w:t></w:r><w:proofErr w:type="gramStart" /><w:r><w:t>made</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:t xml:space="preserve"> they have been very popular. The original Mario Party sold 2.7 million copies. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">As for the potential social impact our product may have a social impact if it takes off because of the genre it is in. Party games </w:t></w:r><w:proofErr w:type="gramStart" /><w:r><w:t>have the opportunity to</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:t xml:space="preserve"> bring people together or </w:t></w:r><w:r><w:lastRenderedPageBreak /><w:t xml:space="preserve">pin them against each other, and for many that is very captivating. If </w:t></w:r><w:proofErr w:type="gramStart" /><w:r><w:t>anything</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:t xml:space="preserve"> our game will inspire people who play it to explore more of the same genre or potentially make their own.  </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Unfortunately, there are lots of similar products to what we are making. Things like Mario Party have been around since before 2000. Even today newer party games like Pummel Party are being made and still seeing success. Our game is going to differ from these but still stick to the traditional board game format that is so popular in party games. With a game like Mario Party that also uses this you need everyone in the same place to play together. Otherwise, you need Nintendo’s subscription service to play together. Our game will not require that. Our game is also going to have one host player that is in control, making it more like a traditional board game in that sense. Everyone else will be contestants making this a more unique concept when it comes to traditional party video games.</w:t></w:r></w:p>'
$sel.InsertXML($xml)
